$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 147 (this pushes old rows 147-162 down to 151-166)
$ws.Range("147:150").Insert()

# Common constant values for these rows
$mercadoId = 11
$mercado = "Vega Monumental Concepción"
$region = "Bíobío"
$codreg = 8
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103001
$categoria = "Cereza"
$kgUnidad = 10

function Set-Row([int]$r, [double]$fecha, [string]$variedad, [string]$calidad, [double]$volumen, [double]$precioMin, [double]$precioMax, [double]$precioProm, [string]$unidad, [string]$origen, [double]$precioKg) {
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

Set-Row 147 44918 "Lapins" "Primera" 100 5000 6000 5500 "`$/bandeja 10 kilos" "Región de O'Higgins" 550
Set-Row 148 44918 "Lapins" "Segunda" 50 4000 4000 4000 "`$/bandeja 10 kilos" "Región de O'Higgins" 400
Set-Row 149 44918 "Sweet Heart" "Primera" 100 5000 6000 5500 "`$/bandeja 10 kilos" "Región de O'Higgins" 550
Set-Row 150 44918 "Sweet Heart" "Segunda" 50 4000 4000 4000 "`$/bandeja 10 kilos" "Región de O'Higgins" 400
